# Apply edits described by the commit diff:
#  - Row 2 data corrections (name typo fix, course hours, passport number)
#  - Selection moved to row 3 (whole row), as if the next empty row were selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix candidate full name typo: "Amirov AKrom" -> "Amirov Akrom"
$ws.Range("A2").Value = "Amirov Akrom"

# Update course/training description: 864 soatlik -> 576 soatlik
$ws.Range("B2").Value = "Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik"

# Update passport number: AA7778899 -> AB9890099
$ws.Range("C2").Value = "AB9890099"

# Move the active selection to the whole of row 3 (next row), mirroring the
# change of the workbook's stored selection from A2:XFD2 to A3:XFD3
$ws.Rows("3").Select()
